# Append new ICD code rows to the "Deaths" sheet, right under the existing
# data block, mirroring how the Python source now determines the next free
# row/column bounds via openpyxl before writing the new records.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Deaths")

$newRows = @(
    @{ Id = 12; Code = "A34" },
    @{ Id = 13; Code = "B31" },
    @{ Id = 14; Code = "G564" }
)

# openpyxl-style bound lookup: last used row is 12, so new rows start at 14
# (row 13 is intentionally skipped, matching the source data export).
$startRow = $ws.UsedRange.Rows.Count + 2

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $entry = $newRows[$i]

    $ws.Cells.Item($row, 1).Value = $entry.Id
    $ws.Cells.Item($row, 2).Value = $entry.Code
}

$ws.Range("C14").Select() | Out-Null
